$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.939.21'
$ws.Range('E2').Value = '  -0.75%  '
$ws.Range('D3').Value = '1.650.64'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E7').Value = '  -1.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3828'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '52.11'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.355'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9997'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08428'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.92'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.089'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.934'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001316'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('D17').Value = '1.647.36'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '94.69'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06960'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.71'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.943'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.50%  '
$ws.Range('E22').Value = '  +0.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.74'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Value = '23.938.47'
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.452'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.978'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '151.85'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.405'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '138.85'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.867'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.524'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.56%  '
$ws.Range('D33').Value = '1.827.41'
$ws.Range('E33').Value = '  +1.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.046'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08069'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02975'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.13%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.668'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.66%  '
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '10.89'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2682'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09099'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.7629'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.47'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.81%  '
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.35'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7021'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.470'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9998'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08304'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '134.60'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.214'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.98%  '
